$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @'
questions = [
    {
        "title": "Company X sold 6,500 units last year. The total number of units sold on the entire market was 357,500 last year.  What was the market share of Company X last year?",
        "ques_type": 2,
        "options": [
            "1.8%",
            "5.5%",
            "0.01%",
            "54.9%"
        ],
        "score": "1.8%"
    },
    {
        "title": "Your client wants to enter the milk market in a new country. One company holds 28% of the market share in that country, while the remaining share of the market is highly fragmented among numerous smaller companies. Your client intends to establish a long-term presence in the market and has sufficient capital to maneuver before and after entry. Which of the following should you advise your client about entering the market?",
        "ques_type": 2,
        "options": [
            "Market entry is not feasible as the client will be forced to compete among the smaller players.",
            "Market entry will not pose difficulties at all as the client can simply drive out smaller players and establish its business subsequently.",
            "The client should consider acquiring one or more existing companies with established infrastructure and presence in the market.",
            "The client should consider entering the country in a product category with greater market fragmentation instead of challenging the dominant company in the milk sector."
        ],
        "score": "The client should consider acquiring one or more existing companies with established infrastructure and presence in the market."
    },
    {
        "title": "You are analyzing the carbonated soft drink industry. Total industry revenue in the United States is expected to hit $144.5 billion this year and grow by 2.65% annually for four years. Per person consumption is 34.3 gallons (130 liters), while per person revenues are $430. Next year, volume growth per person in carbonated soft drinks should be just over 1%, and in four years, the total volume is expected to be just over 10.5 million gallons (40 million liters).  Which of these insights is best supported?",
        "ques_type": 2,
        "options": [
            "The United States consumes the most carbonated soft drinks per capita worldwide.",
            "Volume growth of carbonated soft drinks in the United States is expected to grow at a compound annual growth rate (CAGR) of 2.65% between 2022-2026.",
            "Carbonated and noncarbonated soft drinks are expected to amount to around 21 gallons/80 liters by 2026.",
            "The unit price (price per gallon/liter) is expected to grow in the period between 2022-2026."
        ],
        "score": "The unit price (price per gallon/liter) is expected to grow in the period between 2022-2026."
    },
    {
        "title": "Some brands establish a strong presence, capture a loyal customer base, and hold a strong No. 1 position in a given market.  True or false: When a company is in this position, it should redirect funds from marketing and innovation toward expansion into unexplored markets to further grow sales.",
        "ques_type": 11,
        "options": [
            "true",
            "false"
        ],
        "score": "False"
    }
]
'@

# Drop the second row (the old shared-string cell at A2) entirely.
$ws.Rows(2).Delete()

# Replace A1's content with the pretty-printed questions text.
$ws.Range("A1").Value = $newText

# Restore the row's natural height (undo any autofit from the multi-line text)
# before dropping the old bold/bordered formatting back to the sheet default.
$ws.Rows(1).AutoFit() | Out-Null
$ws.Range("A1").ClearFormats()
